$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure every updated cell keeps plain-text storage (matches original inlineStr cells)
# instead of Excel auto-coercing numeric-looking / percent-looking text into numbers.
$cellRefs = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","D7","E7","D8","E8","B9","C9","D9","E9","B10","C10","D10","E10","B11","C11","D11","E11","B12","C12","D12","E12","B13","C13","D13","E13","B14","C14","D14","E14","B15","C15","D15","E15","D16","E16","E17","E18","D20","E20","D21","E21","D22","E22","E23","D24","D25","E25","D26","E26","E27","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","E46","E47","D48","E48","D49","E49","D50","E50")
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values cell by cell
$ws.Range("D2").Value = "257.28"
$ws.Range("E2").Value = "0.07%"
$ws.Range("D3").Value = "27.06"
$ws.Range("E3").Value = "-1.59%"
$ws.Range("D4").Value = "4.693"
$ws.Range("E4").Value = "-9.91%"
$ws.Range("D5").Value = "0.05880"
$ws.Range("E5").Value = "-0.54%"
$ws.Range("D6").Value = "6.646"
$ws.Range("D7").Value = "0.8579"
$ws.Range("E7").Value = "-1.11%"
$ws.Range("D8").Value = "0.9605"
$ws.Range("E8").Value = "-5.28%"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "0.01051"
$ws.Range("E9").Value = "1,639.88%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1410"
$ws.Range("E10").Value = "-0.59%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.03954"
$ws.Range("E11").Value = "10.46%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07095"
$ws.Range("E12").Value = "-1.27%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03183"
$ws.Range("E13").Value = "0.21%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09172"
$ws.Range("E14").Value = "-0.64%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001538"
$ws.Range("E15").Value = "-0.64%"
$ws.Range("D16").Value = "0.006206"
$ws.Range("E16").Value = "5.74%"
$ws.Range("E17").Value = "1.14%"
$ws.Range("E18").Value = "-1.95%"
$ws.Range("D20").Value = "0.3080"
$ws.Range("E20").Value = "-2.25%"
$ws.Range("D21").Value = "0.1292"
$ws.Range("E21").Value = "-1.09%"
$ws.Range("D22").Value = "3.853"
$ws.Range("E22").Value = "9.48%"
$ws.Range("E23").Value = "1.06%"
$ws.Range("D24").Value = "0.001222"
$ws.Range("D25").Value = "0.004296"
$ws.Range("E25").Value = "-4.90%"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").Value = "0.07%"
$ws.Range("E27").Value = "0.01%"
$ws.Range("D40").Value = "0.03830"
$ws.Range("E40").Value = "0.36%"
$ws.Range("D41").Value = "0.006223"
$ws.Range("E41").Value = "11.15%"
$ws.Range("E42").Value = "0.01%"
$ws.Range("D43").Value = "0.001901"
$ws.Range("E43").Value = "0.06%"
$ws.Range("D44").Value = "0.01144"
$ws.Range("E44").Value = "6.81%"
$ws.Range("D45").Value = "0.00005459"
$ws.Range("E45").Value = "0.60%"
$ws.Range("E46").Value = "0.06%"
$ws.Range("E47").Value = "-44.97%"
$ws.Range("D48").Value = "0.1536"
$ws.Range("E48").Value = "6,954.26%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.06%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.06%"
